# Qn2 complete and format for questions
#
# Adds a third slide to the deck (a copy of the "Qn1" config-graph slide,
# repurposed for Qn2): duplicates slide 2, drops the two "1/3" path
# fraction labels that don't apply to the new path, relabels the
# top-left box "Q1a" -> "Q2a", and turns the former "1/3" label over the
# first oval into a plain "1" (S1 is now reached with probability 1),
# repositioning it to sit directly above the oval.

$p = $ppt.ActivePresentation

# Duplicate slide 2 ("Qn1 config graph") - this becomes the new slide 3
# and keeps the same shapes/positions/connectors as the source slide.
$newSlide = $p.Slides.Item(2).Duplicate()
$s3 = $p.Slides.Item(3)

# Remove the two "1/3" labels that belonged to the third (now unused)
# outgoing path on the duplicated slide.
$s3.Shapes.Item("TextBox 6").Delete()
$s3.Shapes.Item("TextBox 18").Delete()

# "Q1a" -> "Q2a" for the new question's label.
$s3.Shapes.Item("TextBox 21").TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1).Text = "Q2a"

# The "1/3" above Oval 1 becomes "1" and is recentred above the oval.
$topLabel = $s3.Shapes.Item("TextBox 71")
$topLabel.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1).Text = "1"
$topLabel.Left = 2575683 / 12700
$topLabel.Top = 492639 / 12700
